# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "28.294.85" or "219.36". Some of
# these look like plain decimals, so force the column to Text format
# before writing them back -- otherwise Excel's COM layer would silently
# reinterpret a value like "219.36" as the number 219.36 instead of
# keeping it as the original text string.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Cell = 'D2'; Value = '28.294.85' },
    @{ Cell = 'E2'; Value = '  +3.96%  ' },
    @{ Cell = 'D3'; Value = '1.730.46' },
    @{ Cell = 'E3'; Value = '  +2.53%  ' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '219.36' },
    @{ Cell = 'E5'; Value = '  +1.46%  ' },
    @{ Cell = 'E6'; Value = '  +0.02%  ' },
    @{ Cell = 'E7'; Value = '  -0.05%  ' },
    @{ Cell = 'E8'; Value = '  +4.52%  ' },
    @{ Cell = 'D9'; Value = '0.267' },
    @{ Cell = 'E9'; Value = '  +2.42%  ' },
    @{ Cell = 'D10'; Value = '0.0636' },
    @{ Cell = 'E10'; Value = '  +1.41%  ' },
    @{ Cell = 'D11'; Value = '0.0896' },
    @{ Cell = 'E11'; Value = '  +0.60%  ' },
    @{ Cell = 'D12'; Value = '1.976.36' },
    @{ Cell = 'E12'; Value = '  +2.66%  ' },
    @{ Cell = 'D13'; Value = '1.729.73' },
    @{ Cell = 'E13'; Value = '  +1.95%  ' },
    @{ Cell = 'D14'; Value = '4.26' },
    @{ Cell = 'E14'; Value = '  +1.61%  ' },
    @{ Cell = 'D15'; Value = '0.563' },
    @{ Cell = 'E15'; Value = '  +2.02%  ' },
    @{ Cell = 'D16'; Value = '67.77' },
    @{ Cell = 'E16'; Value = '  +0.69%  ' },
    @{ Cell = 'D17'; Value = '28.294.95' },
    @{ Cell = 'E17'; Value = '  +3.98%  ' },
    @{ Cell = 'D18'; Value = '247.49' },
    @{ Cell = 'E18'; Value = '  +3.97%  ' },
    @{ Cell = 'D19'; Value = '0.0₃0754' },
    @{ Cell = 'E19'; Value = '  +1.09%  ' },
    @{ Cell = 'D20'; Value = '7.94' },
    @{ Cell = 'E20'; Value = '  -2.75%  ' },
    @{ Cell = 'E21'; Value = '  -0.21%  ' },
    @{ Cell = 'D22'; Value = '4.65' },
    @{ Cell = 'E22'; Value = '  +1.59%  ' },
    @{ Cell = 'D23'; Value = '9.69' },
    @{ Cell = 'E23'; Value = '  +0.28%  ' },
    @{ Cell = 'E24'; Value = '  -0.52%  ' },
    @{ Cell = 'D25'; Value = '149.50' },
    @{ Cell = 'E25'; Value = '  +0.78%  ' },
    @{ Cell = 'D26'; Value = '7.50' },
    @{ Cell = 'E26'; Value = '  +2.55%  ' },
    @{ Cell = 'D27'; Value = '16.74' },
    @{ Cell = 'E27'; Value = '  +1.37%  ' },
    @{ Cell = 'E28'; Value = '  +0.56%  ' },
    @{ Cell = 'E29'; Value = '  +0.12%  ' },
    @{ Cell = 'E30'; Value = '  +2.63%  ' },
    @{ Cell = 'E31'; Value = '  +2.75%  ' },
    @{ Cell = 'D32'; Value = '3.42' },
    @{ Cell = 'E32'; Value = '  +0.57%  ' },
    @{ Cell = 'D33'; Value = '3.27' },
    @{ Cell = 'E33'; Value = '  +0.92%  ' },
    @{ Cell = 'D34'; Value = '1.487.45' },
    @{ Cell = 'E34'; Value = '  -5.79%  ' },
    @{ Cell = 'E35'; Value = '  -2.25%  ' },
    @{ Cell = 'D36'; Value = '0.977' },
    @{ Cell = 'E36'; Value = '  +1.86%  ' },
    @{ Cell = 'E37'; Value = '  -0.29%  ' },
    @{ Cell = 'E38'; Value = '  +0.60%  ' },
    @{ Cell = 'E39'; Value = '  +1.01%  ' },
    @{ Cell = 'E40'; Value = '  +0.30%  ' },
    @{ Cell = 'D41'; Value = '70.04' },
    @{ Cell = 'E41'; Value = '  +0.59%  ' },
    @{ Cell = 'E42'; Value = '  -0.10%  ' },
    @{ Cell = 'D43'; Value = '5.64' },
    @{ Cell = 'E43'; Value = '  -1.65%  ' },
    @{ Cell = 'D44'; Value = '1.879.29' },
    @{ Cell = 'E44'; Value = '  +2.39%  ' },
    @{ Cell = 'E45'; Value = '  +1.13%  ' },
    @{ Cell = 'D46'; Value = '0.797' },
    @{ Cell = 'E46'; Value = '  +1.19%  ' },
    @{ Cell = 'E47'; Value = '  +7.09%  ' },
    @{ Cell = 'E48'; Value = '  +4.27%  ' },
    @{ Cell = 'D49'; Value = '90.43' },
    @{ Cell = 'E49'; Value = '  -1.01%  ' },
    @{ Cell = 'D50'; Value = '8.18' },
    @{ Cell = 'E50'; Value = '  -0.57%  ' },
    @{ Cell = 'E51'; Value = '  -1.03%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
